# Update column G ("K" - strikeouts) values for rows 2-18 in Sheet1
# per commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 2
    8  = 3
    9  = 2
    10 = 3
    11 = 3
    12 = 3
    13 = 1
    14 = 3
    15 = 2
    16 = 2
    17 = 1
    18 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
